$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.090.60"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.58"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("E4").Value = "  +0.69%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.29"
$ws.Range("E5").Value = "  -0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("E7").Value = "  +0.69%  "

$ws.Range("E8").Value = "  -2.49%  "

$ws.Range("E9").Value = "  -2.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.44"
$ws.Range("E10").Value = "  -6.10%  "

$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.657.39"
$ws.Range("E12").Value = "  +0.06%  "

$ws.Range("E13").Value = "  -1.71%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.529"
$ws.Range("E14").Value = "  -2.92%  "

$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₃0749"
$ws.Range("E15").Value = "  -1.96%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.089.90"
$ws.Range("E16").Value = "  +0.09%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.23"
$ws.Range("E17").Value = "  -1.21%  "

$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "189.85"
$ws.Range("E19").Value = "  -1.64%  "

$ws.Range("E20").Value = "  -2.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.52"
$ws.Range("E21").Value = "  -4.12%  "

$ws.Range("E22").Value = "  -2.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "144.23"
$ws.Range("E23").Value = "  +0.59%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("E25").Value = "  +0.72%  "

$ws.Range("E26").Value = "  -1.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.77"
$ws.Range("E27").Value = "  -1.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.21"
$ws.Range("E28").Value = "  -2.40%  "

$ws.Range("E29").Value = "  -0.69%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0484"
$ws.Range("E30").Value = "  -3.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.17"
$ws.Range("E32").Value = "  -4.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.44"
$ws.Range("E33").Value = "  -0.08%  "

$ws.Range("E34").Value = "  -1.97%  "

$ws.Range("E35").Value = "  -2.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.124.31"
$ws.Range("E36").Value = "  -0.99%  "

$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.522"
$ws.Range("E38").Value = "  -4.39%  "

$ws.Range("E39").Value = "  -1.80%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.85"
$ws.Range("E40").Value = "  -0.41%  "

$ws.Range("E41").Value = "  -1.51%  "

$ws.Range("E42").Value = "  -3.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₆0114"
$ws.Range("E43").Value = "  -0.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "55.14"
$ws.Range("E44").Value = "  -2.72%  "

$ws.Range("E45").Value = "  -1.95%  "

$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.59"
$ws.Range("E48").Value = "  -1.57%  "

$ws.Range("E49").Value = "  +0.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0929"
$ws.Range("E50").Value = "  -3.54%  "

$ws.Range("E51").Value = "  -1.15%  "
